$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 456.66666
$ws.Range("I12").Value2 = 437.5
$ws.Range("K12").Value2 = 437.5
$ws.Range("M12").Value2 = -267.5
$ws.Range("H58").Value2 = 1674.2727
$ws.Range("I58").Value2 = 123.4
$ws.Range("J58").Value2 = 2966.6667
$ws.Range("K58").Value2 = 370.2
$ws.Range("L58").Value2 = 8900.000100000001
$ws.Range("M58").Value2 = -220.2
$ws.Range("N58").Value2 = -9200.000100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 99.333336
$ws.Range("I4").Value2 = 99.333336
$ws.Range("K4").Value2 = 99.333336
$ws.Range("M4").Value2 = 16.666664
$ws.Range("H5").Value2 = 100
$ws.Range("I5").Value2 = 100
$ws.Range("K5").Value2 = 100
$ws.Range("M5").Value2 = 12
$ws.Range("H38").Value2 = 3229.5
$ws.Range("I38").Value2 = 3229.5
$ws.Range("K38").Value2 = 3229.5
$ws.Range("M38").Value2 = -2762.5
$ws.Range("H45").Value2 = 3423.2307
$ws.Range("I45").Value2 = 2285.5715
$ws.Range("J45").Value2 = 4750.5
$ws.Range("K45").Value2 = 2285.5715
$ws.Range("L45").Value2 = 4750.5
$ws.Range("M45").Value2 = -1908.5715
$ws.Range("N45").Value2 = -5504.5
$ws.Range("H96").Value2 = 4032818.2
$ws.Range("J96").Value2 = 4032818.2
$ws.Range("L96").Value2 = 4032818.2
$ws.Range("N96").Value2 = -4038310.2
$ws.Range("H97").Value2 = 1927.5
$ws.Range("I97").Value2 = 1450.3334
$ws.Range("K97").Value2 = 1450.3334
$ws.Range("M97").Value2 = -954.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 100
$ws.Range("I4").Value2 = 100
$ws.Range("K4").Value2 = 100
$ws.Range("M4").Value2 = 15
$ws.Range("H94").Value2 = 735.6923
$ws.Range("I94").Value2 = 679.3333
$ws.Range("K94").Value2 = 679.3333
$ws.Range("M94").Value2 = -228.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value2 = 109.666664
$ws.Range("I7").Value2 = 123.833336
$ws.Range("J7").Value2 = 81.333336
$ws.Range("K7").Value2 = 123.833336
$ws.Range("L7").Value2 = 81.333336
$ws.Range("M7").Value2 = -10.833336
$ws.Range("N7").Value2 = -307.333336
$ws.Range("H32").Value2 = 1541.1666
$ws.Range("I32").Value2 = 1399.4
$ws.Range("K32").Value2 = 1399.4
$ws.Range("M32").Value2 = -1083.4
$ws.Range("H33").Value2 = 1200
$ws.Range("I33").Value2 = 1200
$ws.Range("K33").Value2 = 1200
$ws.Range("M33").Value2 = -821
$ws.Range("H36").Value2 = 5298.5
$ws.Range("I36").Value2 = 5298.5
$ws.Range("K36").Value2 = 5298.5
$ws.Range("M36").Value2 = -4910.5
$ws.Range("H40").Value2 = 5298.5
$ws.Range("I40").Value2 = 5298.5
$ws.Range("K40").Value2 = 5298.5
$ws.Range("M40").Value2 = -5138.5
$ws.Range("H44").Value2 = 30000
$ws.Range("I44").Value2 = 30000
$ws.Range("K44").Value2 = 30000
$ws.Range("M44").Value2 = -29558
$ws.Range("H55").Value2 = 0
$ws.Range("I55").Value2 = 0
$ws.Range("K55").Value2 = 0
$ws.Range("M55").ClearContents()
$ws.Range("H99").Value2 = 6749.75
$ws.Range("I99").Value2 = 7500
$ws.Range("J99").Value2 = 5999.5
$ws.Range("K99").Value2 = 7500
$ws.Range("L99").Value2 = 5999.5
$ws.Range("M99").Value2 = -6002
$ws.Range("N99").Value2 = -8995.5
$ws.Range("H126").Value2 = 6749.75
$ws.Range("I126").Value2 = 7500
$ws.Range("J126").Value2 = 5999.5
$ws.Range("K126").Value2 = 22500
$ws.Range("L126").Value2 = 17998.5
$ws.Range("M126").Value2 = -20030
$ws.Range("N126").Value2 = -22938.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value2 = 354.36365
$ws.Range("I14").Value2 = 354.36365
$ws.Range("K14").Value2 = 1063.09095
$ws.Range("M14").Value2 = -890.09095
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value2 = 0
$ws.Range("H17").Value2 = 430.85715
$ws.Range("I17").Value2 = 80.69231000000001
$ws.Range("J17").Value2 = 999.875
$ws.Range("K17").Value2 = 242.07693
$ws.Range("L17").Value2 = 2999.625
$ws.Range("M17").Value2 = -73.07693
$ws.Range("N17").Value2 = -3337.625
$ws.Range("H41").Value2 = 200
$ws.Range("I41").Value2 = 0
$ws.Range("K41").Value2 = 0
$ws.Range("M41").ClearContents()
$ws.Range("H68").Value2 = 810
$ws.Range("J68").Value2 = 767.5
$ws.Range("L68").Value2 = 2302.5
$ws.Range("N68").Value2 = -3924.5
$ws.Range("H71").Value2 = 810
$ws.Range("J71").Value2 = 767.5
$ws.Range("L71").Value2 = 6907.5
$ws.Range("N71").Value2 = -15019.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 100032.4
$ws.Range("I2").Value2 = 250013.75
$ws.Range("J2").Value2 = 44.833332
$ws.Range("K2").Value2 = 250013.75
$ws.Range("L2").Value2 = 44.833332
$ws.Range("M2").Value2 = -249900.75
$ws.Range("N2").Value2 = -270.833332
$ws.Range("H99").Value2 = 8096.3335
$ws.Range("I99").Value2 = 4144.5
$ws.Range("K99").Value2 = 4144.5
$ws.Range("M99").Value2 = -1898.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 1242.2858
$ws.Range("I55").Value2 = 1594.3334
$ws.Range("J55").Value2 = 978.25
$ws.Range("K55").Value2 = 1594.3334
$ws.Range("L55").Value2 = 978.25
$ws.Range("M55").Value2 = -1421.3334
$ws.Range("N55").Value2 = -1324.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value2 = 16428571
$ws.Range("J5").Value2 = 20000000
$ws.Range("L5").Value2 = 20000000
$ws.Range("N5").Value2 = -20000224
$ws.Range("H107").Value2 = 780
$ws.Range("I107").Value2 = 638.6
$ws.Range("J107").Value2 = 956.75
$ws.Range("K107").Value2 = 1915.8
$ws.Range("L107").Value2 = 2870.25
$ws.Range("M107").Value2 = 4.199999999999818
$ws.Range("N107").Value2 = -6710.25
$ws.Range("H126").Value2 = 5500
$ws.Range("I126").Value2 = 3500
$ws.Range("K126").Value2 = 10500
$ws.Range("M126").Value2 = -8030
$ws.Range("H132").Value2 = 4618.6665
$ws.Range("I132").Value2 = 5000
$ws.Range("J132").Value2 = 4571
$ws.Range("K132").Value2 = 15000
$ws.Range("L132").Value2 = 13713
$ws.Range("M132").Value2 = -12470
$ws.Range("N132").Value2 = -18773
